$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '$2b$15$gOcnXVvpIuSEpDX5aajxQ.'
$ws.Range("C2").Value = '$2b$15$gOcnXVvpIuSEpDX5aajxQ.UdoBc43HOmGik3lv.MQjVNC1.Dxw/Qy'
$ws.Range("D2").Value = 'MBOHIOSMBXXFAGEQOYTCMBBOJRSELO'
